$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimize_config")

# Update template/target string values. D2 and D4 both originally shared the
# "ferrite" string, so both must be updated to keep sharing the new string.
$ws.Range("D2").Value = "ferrite_n20"
$ws.Range("D4").Value = "ferrite_n20"
$ws.Range("G2").Value = "0p1"
$ws.Range("D5").Value = "ferrite_n100"

# Adjust column F width to match new (shorter) best-fit content. The stored
# width is quantized to a pixel grid by the host, so feed it the value that
# rounds closest to the target 20.719285714285714 character-width.
$ws.Columns.Item(6).ColumnWidth = 19.75
